# Add the TP11 row (row 15) to Sheet1, matching the pattern of the
# preceding "TPx" rows: col A = "TPx" label, col B = long comment text,
# col C = grade. Row 13 (TP10) uses the style trio we want to reuse
# (A: left/top on shaded fill, B: wrap+quote-prefix on shaded fill,
# C: centered percent-style on shaded fill), so we set the new cells'
# values first and then copy just the *formatting* from row 13 onto
# row 15 - this keeps the existing style indices (no new styles are
# minted) and avoids disturbing the SUM formula's dependency tracking
# (which otherwise goes stale if the format copy happens before the
# values are written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "TP11"

$comment = "Visualização arquitetural (dois apresentaram)`n- antlr4, json, java, mayavi, networkx`n- análise coleta tudo em json e plotar`n- grafo, mas com cor para cada componente.`n- Fizeram algo buscando um code city.`n- Atividade desempenhada de forma como esperada."
$ws.Range("B15").Value = $comment

$ws.Range("C15").Value = 7.5

# Pick up the A/B/C formatting (fills, alignment, wrap, quote-prefix,
# number format) from the row above (TP10) so row 15 looks like its
# siblings.
$ws.Range("A13:C13").Copy()
$ws.Range("A15").PasteSpecial(-4122)

# Match the row height from the target workbook.
$ws.Rows.Item(15).RowHeight = 102

# Selection moved to E13 (and the view scrolled down) after the edit.
$ws.Range("E13").Select()
